$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value2 = 1847.2222
$ws.Range("J40").Value2 = 2089.5715
$ws.Range("L40").Value2 = 2089.5715
$ws.Range("N40").Value2 = -2439.5715

$ws.Range("H92").Value2 = 25006296
$ws.Range("I92").Value2 = 31257624
$ws.Range("J92").Value2 = 977.5
$ws.Range("K92").Value2 = 31257624
$ws.Range("L92").Value2 = 977.5
$ws.Range("M92").Value2 = -31256376
$ws.Range("N92").Value2 = -3473.5

$ws.Range("H113").Value2 = 2264.9
$ws.Range("I113").Value2 = 1925
$ws.Range("J113").Value2 = 2491.5
$ws.Range("K113").Value2 = 1925
$ws.Range("L113").Value2 = 2491.5
$ws.Range("M113").Value2 = 1329
$ws.Range("N113").Value2 = -8999.5

$ws.Range("H134").Value2 = 0
$ws.Range("J134").Value2 = 0
$ws.Range("L134").Value2 = 0
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value2 = 0
$ws.Range("J136").Value2 = 0
$ws.Range("L136").Value2 = 0
$ws.Range("N136").ClearContents()

$ws.Range("H137").Value2 = 20925.42
$ws.Range("I137").Value2 = 25220.268
$ws.Range("J137").Value2 = 1360
$ws.Range("K137").Value2 = 75660.804
$ws.Range("L137").Value2 = 4080
$ws.Range("M137").Value2 = -73110.804
$ws.Range("N137").Value2 = -9180

$ws.Range("H139").Value2 = 0
$ws.Range("J139").Value2 = 0
$ws.Range("L139").Value2 = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 1531.7959
$ws.Range("I61").Value2 = 1560.409
$ws.Range("J61").Value2 = 1280
$ws.Range("K61").Value2 = 1560.409
$ws.Range("L61").Value2 = 1280
$ws.Range("M61").Value2 = -1348.409
$ws.Range("N61").Value2 = -1704

$ws.Range("H122").Value2 = 1313.091
$ws.Range("I122").Value2 = 1222.1666
$ws.Range("J122").Value2 = 1422.2
$ws.Range("K122").Value2 = 3666.4998
$ws.Range("L122").Value2 = 4266.6
$ws.Range("M122").Value2 = -1216.4998
$ws.Range("N122").Value2 = -9166.6

$ws.Range("H136").Value2 = 1531.7959
$ws.Range("I136").Value2 = 1560.409
$ws.Range("J136").Value2 = 1280
$ws.Range("K136").Value2 = 4681.227000000001
$ws.Range("L136").Value2 = 3840
$ws.Range("M136").Value2 = -2131.227000000001
$ws.Range("N136").Value2 = -8940

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 1754.5454
$ws.Range("I20").Value2 = 1780
$ws.Range("J20").Value2 = 1500
$ws.Range("K20").Value2 = 1780
$ws.Range("L20").Value2 = 1500
$ws.Range("M20").Value2 = -1533
$ws.Range("N20").Value2 = -1994

$ws.Range("H22").Value2 = 1149.7142
$ws.Range("I22").Value2 = 2775.25
$ws.Range("J22").Value2 = 499.5
$ws.Range("K22").Value2 = 2775.25
$ws.Range("L22").Value2 = 499.5
$ws.Range("M22").Value2 = -2602.25
$ws.Range("N22").Value2 = -845.5

$ws.Range("H105").Value2 = 1301.3334
$ws.Range("I105").Value2 = 0
$ws.Range("J105").Value2 = 1301.3334
$ws.Range("K105").Value2 = 0
$ws.Range("L105").Value2 = 1301.3334
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value2 = -4795.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 7790.6113
$ws.Range("I31").Value2 = 6503.4863
$ws.Range("K31").Value2 = 6503.4863
$ws.Range("M31").Value2 = -6208.4863

$ws.Range("H34").Value2 = 7790.6113
$ws.Range("I34").Value2 = 6503.4863
$ws.Range("K34").Value2 = 6503.4863
$ws.Range("M34").Value2 = -6301.4863

$ws.Range("H105").Value2 = 752.25
$ws.Range("I105").Value2 = 666
$ws.Range("J105").Value2 = 1011
$ws.Range("K105").Value2 = 666
$ws.Range("L105").Value2 = 1011
$ws.Range("M105").Value2 = 1081
$ws.Range("N105").Value2 = -4505

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value2 = 260.2
$ws.Range("I13").Value2 = 260.2
$ws.Range("K13").Value2 = 780.5999999999999
$ws.Range("M13").Value2 = -612.5999999999999

$ws.Range("H64").Value2 = 2528604.8
$ws.Range("I64").Value2 = 1187.3334
$ws.Range("J64").Value2 = 5056022.5
$ws.Range("K64").Value2 = 3562.0002
$ws.Range("L64").Value2 = 15168067.5
$ws.Range("M64").Value2 = -3292.0002
$ws.Range("N64").Value2 = -15168607.5

$ws.Range("H67").Value2 = 2528604.8
$ws.Range("I67").Value2 = 1187.3334
$ws.Range("J67").Value2 = 5056022.5
$ws.Range("K67").Value2 = 3562.0002
$ws.Range("L67").Value2 = 15168067.5
$ws.Range("M67").Value2 = -2626.0002
$ws.Range("N67").Value2 = -15169939.5

$ws.Range("H87").Value2 = 57750.684
$ws.Range("J87").Value2 = 57750.684
$ws.Range("L87").Value2 = 173252.052
$ws.Range("N87").Value2 = -175748.052

$ws.Range("H90").Value2 = 57750.684
$ws.Range("J90").Value2 = 57750.684
$ws.Range("L90").Value2 = 519756.156
$ws.Range("N90").Value2 = -532236.156

$ws.Range("H122").Value2 = 467
$ws.Range("I122").Value2 = 294
$ws.Range("K122").Value2 = 2646
$ws.Range("M122").Value2 = -196

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value2 = 12819.125
$ws.Range("I43").Value2 = 678.75
$ws.Range("J43").Value2 = 24959.5
$ws.Range("K43").Value2 = 678.75
$ws.Range("L43").Value2 = 24959.5
$ws.Range("M43").Value2 = -527.75
$ws.Range("N43").Value2 = -25261.5

$ws.Range("H46").Value2 = 9446
$ws.Range("I46").Value2 = 0
$ws.Range("J46").Value2 = 9446
$ws.Range("K46").Value2 = 0
$ws.Range("L46").Value2 = 9446
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value2 = -9758

$ws.Range("H80").Value2 = 6077.1377
$ws.Range("I80").Value2 = 2826.818
$ws.Range("J80").Value2 = 8063.4443
$ws.Range("K80").Value2 = 2826.818
$ws.Range("L80").Value2 = 8063.4443
$ws.Range("M80").Value2 = -1828.818
$ws.Range("N80").Value2 = -10059.4443

$ws.Range("H83").Value2 = 6077.1377
$ws.Range("I83").Value2 = 2826.818
$ws.Range("J83").Value2 = 8063.4443
$ws.Range("K83").Value2 = 14134.09
$ws.Range("L83").Value2 = 40317.2215
$ws.Range("M83").Value2 = -9142.09
$ws.Range("N83").Value2 = -50301.2215

$ws.Range("H122").Value2 = 2416.375
$ws.Range("I122").Value2 = 2883
$ws.Range("J122").Value2 = 1389.8
$ws.Range("K122").Value2 = 8649
$ws.Range("L122").Value2 = 4169.4
$ws.Range("M122").Value2 = -6199
$ws.Range("N122").Value2 = -9069.4

$ws.Range("H132").Value2 = 25328.453
$ws.Range("I132").Value2 = 1479.1428
$ws.Range("J132").Value2 = 49177.76
$ws.Range("K132").Value2 = 4437.428400000001
$ws.Range("L132").Value2 = 147533.28
$ws.Range("M132").Value2 = -1907.428400000001
$ws.Range("N132").Value2 = -152593.28

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value2 = 2573.9412
$ws.Range("J122").Value2 = 2292.2222
$ws.Range("L122").Value2 = 6876.6666
$ws.Range("N122").Value2 = -11776.6666

$ws.Range("H136").Value2 = 258076.16
$ws.Range("I136").Value2 = 436108.53
$ws.Range("K136").Value2 = 1308325.59
$ws.Range("M136").Value2 = -1305775.59
